$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (with thousand separators as dots, trailing zeros, etc.) are preserved
# exactly as authored rather than being auto-coerced to Double values.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "26.278.87"
$ws.Range("E2").Value = "  +0.09%  "

# Row 3
$ws.Range("D3").Value = "1.681.04"
$ws.Range("E3").Value = "  +0.40%  "

# Row 4
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.25%  "

# Row 5
$ws.Range("D5").Value = "218.19"
$ws.Range("E5").Value = "  +0.04%  "

# Row 6
$ws.Range("D6").Value = "0.5270"
$ws.Range("E6").Value = "  +2.61%  "

# Row 7
$ws.Range("E7").Value = "  +0.24%  "

# Row 8
$ws.Range("D8").Value = "0.2700"
$ws.Range("E8").Value = "  +1.43%  "

# Row 9
$ws.Range("D9").Value = "0.06422"
$ws.Range("E9").Value = "  +0.20%  "

# Row 10
$ws.Range("D10").Value = "22.06"
$ws.Range("E10").Value = "  +2.31%  "

# Row 11
$ws.Range("D11").Value = "0.07490"
$ws.Range("E11").Value = "  +1.58%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.703.01"
$ws.Range("E12").Value = "  +1.70%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.549"
$ws.Range("E13").Value = "  -0.13%  "

# Row 14
$ws.Range("E14").Value = "  -0.31%  "

# Row 15
$ws.Range("D15").Value = "0.000008479"
$ws.Range("E15").Value = "  -2.49%  "

# Row 16
$ws.Range("D16").Value = "64.30"
$ws.Range("E16").Value = "  -0.93%  "

# Row 17
$ws.Range("D17").Value = "26.323.16"
$ws.Range("E17").Value = "  -0.01%  "

# Row 18
$ws.Range("D18").Value = "4.925"
$ws.Range("E18").Value = "  -0.68%  "

# Row 19
$ws.Range("D19").Value = "1.008"
$ws.Range("E19").Value = "  +0.19%  "

# Row 20
$ws.Range("D20").Value = "10.88"
$ws.Range("E20").Value = "  +0.28%  "

# Row 21
$ws.Range("D21").Value = "188.87"
$ws.Range("E21").Value = "  -0.67%  "

# Row 22
$ws.Range("D22").Value = "6.198"
$ws.Range("E22").Value = "  -0.45%  "

# Row 23
$ws.Range("D23").Value = "1.010"
$ws.Range("E23").Value = "  +0.29%  "

# Row 24
$ws.Range("D24").Value = "144.38"
$ws.Range("E24").Value = "  +0.03%  "

# Row 25
$ws.Range("D25").Value = "7.721"
$ws.Range("E25").Value = "  +1.16%  "

# Row 26
$ws.Range("D26").Value = "0.1239"
$ws.Range("E26").Value = "  +4.63%  "

# Row 27
$ws.Range("D27").Value = "15.79"
$ws.Range("E27").Value = "  +0.97%  "

# Row 28
$ws.Range("D28").Value = "0.06610"
$ws.Range("E28").Value = "  +11.38%  "

# Row 29
$ws.Range("E29").Value = "  +4.92%  "

# Row 30
$ws.Range("E30").Value = "  +0.42%  "

# Row 31
$ws.Range("D31").Value = "3.581"
$ws.Range("E31").Value = "  +1.55%  "

# Row 32
$ws.Range("D32").Value = "3.562"
$ws.Range("E32").Value = "  +0.72%  "

# Row 33
$ws.Range("D33").Value = "1.663"
$ws.Range("E33").Value = "  +1.44%  "

# Row 34
$ws.Range("D34").Value = "1.027"
$ws.Range("E34").Value = "  +1.14%  "

# Row 35
$ws.Range("D35").Value = "0.6200"
$ws.Range("E35").Value = "  +3.00%  "

# Row 36
$ws.Range("D36").Value = "2.398"
$ws.Range("E36").Value = "  +1.28%  "

# Row 37
$ws.Range("D37").Value = "2.698"
$ws.Range("E37").Value = "  +1.75%  "

# Row 38
$ws.Range("D38").Value = "6.389"
$ws.Range("E38").Value = "  +5.38%  "

# Row 39
$ws.Range("D39").Value = "1.109.63"
$ws.Range("E39").Value = "  +2.95%  "

# Row 40
$ws.Range("E40").Value = "  +0.19%  "

# Row 41
$ws.Range("D41").Value = "0.8766"
$ws.Range("E41").Value = "  +0.87%  "

# Row 42
$ws.Range("E42").Value = "  +0.41%  "

# Row 43
$ws.Range("D43").Value = "100.64"

# Row 44
$ws.Range("D44").Value = "1.830.43"
$ws.Range("E44").Value = "  +0.45%  "

# Row 45
$ws.Range("D45").Value = "0.00000000115"
$ws.Range("E45").Value = "  +2.46%  "

# Row 46
$ws.Range("D46").Value = "56.83"
$ws.Range("E46").Value = "  +1.43%  "

# Row 47
$ws.Range("D47").Value = "8.148"
$ws.Range("E47").Value = "  +1.04%  "

# Row 48
$ws.Range("D48").Value = "1.003"
$ws.Range("E48").Value = "  -0.51%  "

# Row 49
$ws.Range("D49").Value = "0.05272"

# Row 50
$ws.Range("D50").Value = "0.4307"
$ws.Range("E50").Value = "  +0.16%  "

# Row 51
$ws.Range("D51").Value = "6.046"
$ws.Range("E51").Value = "  +2.77%  "
